$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1803.4
$ws.Range("I38").Value = 106.8
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 320.4
$ws.Range("L38").Value = 10500
$ws.Range("M38").Value = 51.60000000000002
$ws.Range("N38").Value = -11244
$ws.Range("H58").Value = 1260.4546
$ws.Range("I58").Value = 266.42856
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 799.28568
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -649.28568
$ws.Range("N58").Value = -9300
$ws.Range("H70").Value = 811.7059
$ws.Range("I70").Value = 709.4524
$ws.Range("J70").Value = 976.88464
$ws.Range("K70").Value = 2128.3572
$ws.Range("L70").Value = 2930.65392
$ws.Range("M70").Value = -1858.3572
$ws.Range("N70").Value = -3470.65392
$ws.Range("H73").Value = 811.7059
$ws.Range("I73").Value = 709.4524
$ws.Range("J73").Value = 976.88464
$ws.Range("K73").Value = 2128.3572
$ws.Range("L73").Value = 2930.65392
$ws.Range("M73").Value = -1192.3572
$ws.Range("N73").Value = -4802.65392

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 39600
$ws.Range("J23").Value = 39600
$ws.Range("L23").Value = 39600
$ws.Range("N23").Value = -40118
$ws.Range("H25").Value = 4903.2
$ws.Range("I25").Value = 1504
$ws.Range("J25").Value = 18500
$ws.Range("K25").Value = 1504
$ws.Range("L25").Value = 18500
$ws.Range("M25").Value = -1102
$ws.Range("N25").Value = -19304
$ws.Range("H27").Value = 20000
$ws.Range("J27").Value = 20000
$ws.Range("L27").Value = 20000
$ws.Range("N27").Value = -20368
$ws.Range("H32").Value = 7280.3184
$ws.Range("I32").Value = 3953.8845
$ws.Range("K32").Value = 3953.8845
$ws.Range("M32").Value = -3666.8845
$ws.Range("H35").Value = 8026
$ws.Range("I35").Value = 2518.5
$ws.Range("J35").Value = 19041
$ws.Range("K35").Value = 2518.5
$ws.Range("L35").Value = 19041
$ws.Range("M35").Value = -2112.5
$ws.Range("N35").Value = -19853
$ws.Range("H61").Value = 1412.25
$ws.Range("I61").Value = 1052.5518
$ws.Range("J61").Value = 2360.5454
$ws.Range("K61").Value = 1052.5518
$ws.Range("L61").Value = 2360.5454
$ws.Range("M61").Value = -840.5518
$ws.Range("N61").Value = -2784.5454
$ws.Range("H122").Value = 1922.585
$ws.Range("I122").Value = 1760.2046
$ws.Range("K122").Value = 5280.6138
$ws.Range("M122").Value = -2830.6138
$ws.Range("H136").Value = 1412.25
$ws.Range("I136").Value = 1052.5518
$ws.Range("J136").Value = 2360.5454
$ws.Range("K136").Value = 3157.6554
$ws.Range("L136").Value = 7081.6362
$ws.Range("M136").Value = -607.6553999999996
$ws.Range("N136").Value = -12181.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 19000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 19000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 19000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -19578
$ws.Range("H36").Value = 3360.4375
$ws.Range("I36").Value = 1063.5385
$ws.Range("J36").Value = 13313.667
$ws.Range("K36").Value = 1063.5385
$ws.Range("L36").Value = 13313.667
$ws.Range("M36").Value = -529.5385000000001
$ws.Range("N36").Value = -14381.667
$ws.Range("H37").Value = 1345.3334
$ws.Range("I37").Value = 1345.3334
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1345.3334
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1208.3334
$ws.Range("N37").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7693889
$ws.Range("I31").Value = 6061576.5
$ws.Range("J31").Value = 10528958
$ws.Range("K31").Value = 6061576.5
$ws.Range("L31").Value = 10528958
$ws.Range("M31").Value = -6061281.5
$ws.Range("N31").Value = -10529548
$ws.Range("H34").Value = 7693889
$ws.Range("I34").Value = 6061576.5
$ws.Range("J34").Value = 10528958
$ws.Range("K34").Value = 6061576.5
$ws.Range("L34").Value = 10528958
$ws.Range("M34").Value = -6061374.5
$ws.Range("N34").Value = -10529362
$ws.Range("H58").Value = 2279.1943
$ws.Range("I58").Value = 546.8570999999999
$ws.Range("J58").Value = 3381.5908
$ws.Range("K58").Value = 546.8570999999999
$ws.Range("L58").Value = 3381.5908
$ws.Range("M58").Value = -343.8570999999999
$ws.Range("N58").Value = -3787.5908
$ws.Range("H132").Value = 2196.45
$ws.Range("I132").Value = 995.1429000000001
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 2985.4287
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -455.4287000000004
$ws.Range("N132").Value = -20058.5
$ws.Range("H133").Value = 25000
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -30060
$ws.Range("H134").Value = 2490.6978
$ws.Range("I134").Value = 2313.4707
$ws.Range("J134").Value = 3160.2222
$ws.Range("K134").Value = 6940.4121
$ws.Range("L134").Value = 9480.6666
$ws.Range("M134").Value = -4405.4121
$ws.Range("N134").Value = -14550.6666
$ws.Range("H135").Value = 37222444
$ws.Range("J135").Value = 37222444
$ws.Range("L135").Value = 37222444
$ws.Range("N135").Value = -37232584
$ws.Range("H136").Value = 2279.1943
$ws.Range("I136").Value = 546.8570999999999
$ws.Range("J136").Value = 3381.5908
$ws.Range("K136").Value = 1640.5713
$ws.Range("L136").Value = 10144.7724
$ws.Range("M136").Value = 909.4287000000002
$ws.Range("N136").Value = -15244.7724
$ws.Range("H137").Value = 16333.333
$ws.Range("I137").Value = 10000
$ws.Range("J137").Value = 29000
$ws.Range("K137").Value = 10000
$ws.Range("L137").Value = 29000
$ws.Range("M137").Value = -4900
$ws.Range("N137").Value = -39200
$ws.Range("H138").Value = 28000
$ws.Range("J138").Value = 28000
$ws.Range("L138").Value = 28000
$ws.Range("N138").Value = -38280
$ws.Range("H140").Value = 41800
$ws.Range("J140").Value = 41800
$ws.Range("L140").Value = 41800
$ws.Range("N140").Value = -52160
$ws.Range("H141").Value = 63500
$ws.Range("J141").Value = 63500
$ws.Range("L141").Value = 63500
$ws.Range("N141").Value = -73860

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1228.5
$ws.Range("I68").Value = 886.75
$ws.Range("J68").Value = 1433.55
$ws.Range("K68").Value = 2660.25
$ws.Range("L68").Value = 4300.65
$ws.Range("M68").Value = -1849.25
$ws.Range("N68").Value = -5922.65
$ws.Range("H71").Value = 1228.5
$ws.Range("I71").Value = 886.75
$ws.Range("J71").Value = 1433.55
$ws.Range("K71").Value = 7980.75
$ws.Range("L71").Value = 12901.95
$ws.Range("M71").Value = -3924.75
$ws.Range("N71").Value = -21013.95
$ws.Range("H131").Value = 839.9231
$ws.Range("I131").Value = 554.0909
$ws.Range("J131").Value = 986.1627999999999
$ws.Range("K131").Value = 1662.2727
$ws.Range("L131").Value = 2958.4884
$ws.Range("M131").Value = 3377.7273
$ws.Range("N131").Value = -13038.4884

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1704.3334
$ws.Range("I41").Value = 1704.3334
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1704.3334
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1349.3334
$ws.Range("N41").ClearContents()
$ws.Range("H113").Value = 1851.15
$ws.Range("I113").Value = 1326.5883
$ws.Range("K113").Value = 1326.5883
$ws.Range("M113").Value = 843.4117000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7819635
$ws.Range("I136").Value = 9623781
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 28871343
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = -28868793
$ws.Range("N136").Value = -10099.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3766.2222
$ws.Range("I136").Value = 912.8
$ws.Range("J136").Value = 18033.334
$ws.Range("K136").Value = 2738.4
$ws.Range("L136").Value = 54100.00199999999
$ws.Range("M136").Value = -188.3999999999996
$ws.Range("N136").Value = -59200.00199999999
